$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 29. This shifts the old row 29
# down to row 30, and the old row 30 down to row 31, preserving their data.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly price record.
$ws.Range("A29").Value = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("D29").Value = 44753
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 100112040
$ws.Range("G29").Value = "Cilantro"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 9000
$ws.Range("L29").Value = 9000
$ws.Range("M29").Value = 9000
$ws.Range("N29").Value = "$/caja 36 atados"
$ws.Range("O29").Value = "Región Metropolitana"
$ws.Range("P29").Value = 250
$ws.Range("Q29").Value = 36
$ws.Range("R29").Value = "Hortaliza"
